$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (row 5)
$ws.Range("A5").Value = 200977
$ws.Range("B5").Value = "Adelola John"
$ws.Range("C5").Value = "adelolajohn@gmail.com"
$ws.Range("D5").Value = "Civil Eng."
$ws.Range("E5").Value = 300
$ws.Range("F5").Value = 564473

# Add hyperlink for the email cell, matching the mailto: pattern used by the
# existing rows, and apply the same built-in "Hyperlink" style used by C2:C4.
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:adelolajohn@gmail.com")
$ws.Range("C5").Style = "Hyperlink"

# Move the active selection to F5 (matches the saved selection in the diff)
$ws.Range("F5").Select()
